$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1661.4286
$ws.Range("I40").Value = 1815
$ws.Range("K40").Value = 1815
$ws.Range("M40").Value = -1640

$ws.Range("H112").Value = 2333
$ws.Range("I112").Value = 2333
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 6999
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("M112").Value = -5891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 6933
$ws.Range("I39").Value = 799
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 799
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -279
$ws.Range("N39").Value = -11040

$ws.Range("H45").Value = 2743.818
$ws.Range("I45").Value = 2568.25
$ws.Range("J45").Value = 4499.5
$ws.Range("K45").Value = 2568.25
$ws.Range("L45").Value = 4499.5
$ws.Range("M45").Value = -2191.25
$ws.Range("N45").Value = -5253.5

$ws.Range("H88").Value = 1483.6471
$ws.Range("I88").Value = 990.8333
$ws.Range("J88").Value = 1752.4546
$ws.Range("K88").Value = 990.8333
$ws.Range("L88").Value = 1752.4546
$ws.Range("M88").Value = -584.8333
$ws.Range("N88").Value = -2564.4546

$ws.Range("H91").Value = 1483.6471
$ws.Range("I91").Value = 990.8333
$ws.Range("J91").Value = 1752.4546
$ws.Range("K91").Value = 990.8333
$ws.Range("L91").Value = 1752.4546
$ws.Range("M91").Value = 413.1667
$ws.Range("N91").Value = -4560.4546

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H95").Value = 32207.5
$ws.Range("J95").Value = 32207.5
$ws.Range("L95").Value = 32207.5
$ws.Range("N95").Value = -37699.5

$ws.Range("H98").Value = 22652
$ws.Range("J98").Value = 22652
$ws.Range("L98").Value = 22652
$ws.Range("N98").Value = -28642

$ws.Range("H102").Value = 3193.818
$ws.Range("I102").Value = 3125.7778
$ws.Range("K102").Value = 3125.7778
$ws.Range("M102").Value = -1503.7778

$ws.Range("H122").Value = 1850
$ws.Range("J122").Value = 1900
$ws.Range("L122").Value = 5700
$ws.Range("N122").Value = -10600

$ws.Range("H132").Value = 6942.8486
$ws.Range("I132").Value = 5900.609
$ws.Range("J132").Value = 9340
$ws.Range("K132").Value = 17701.827
$ws.Range("L132").Value = 28020
$ws.Range("M132").Value = -15171.827
$ws.Range("N132").Value = -33080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4983.9287
$ws.Range("I20").Value = 2210.5557
$ws.Range("K20").Value = 2210.5557
$ws.Range("M20").Value = -1963.5557

$ws.Range("H99").Value = 1107
$ws.Range("I99").Value = 717.36365
$ws.Range("K99").Value = 717.36365
$ws.Range("M99").Value = 780.63635

$ws.Range("H105").Value = 3599.375
$ws.Range("I105").Value = 3299.1667
$ws.Range("K105").Value = 3299.1667
$ws.Range("M105").Value = -1552.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 13205.8
$ws.Range("J28").Value = 13205.8
$ws.Range("L28").Value = 13205.8
$ws.Range("N28").Value = -13695.8

$ws.Range("H43").Value = 29966.334
$ws.Range("J43").Value = 29966.334
$ws.Range("L43").Value = 29966.334
$ws.Range("N43").Value = -30334.334

$ws.Range("H62").Value = 3246.5
$ws.Range("I62").Value = 2494.75
$ws.Range("K62").Value = 2494.75
$ws.Range("M62").Value = -1870.75

$ws.Range("H63").Value = 100271
$ws.Range("J63").Value = 100271
$ws.Range("L63").Value = 100271
$ws.Range("N63").Value = -101643

$ws.Range("H65").Value = 3246.5
$ws.Range("I65").Value = 2494.75
$ws.Range("K65").Value = 12473.75
$ws.Range("M65").Value = -9353.75

$ws.Range("H66").Value = 100271
$ws.Range("J66").Value = 100271
$ws.Range("L66").Value = 300813
$ws.Range("N66").Value = -307677

$ws.Range("H101").Value = 29966.334
$ws.Range("J101").Value = 29966.334
$ws.Range("L101").Value = 29966.334
$ws.Range("N101").Value = -36456.334

$ws.Range("H141").Value = 555523.75
$ws.Range("J141").Value = 555523.75
$ws.Range("L141").Value = 555523.75
$ws.Range("N141").Value = -565883.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 871.75
$ws.Range("I51").Value = 495.66666
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1486.99998
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = -1026.99998
$ws.Range("N51").Value = -6920

$ws.Range("H56").Value = 9147.429
$ws.Range("I56").Value = 9147.429
$ws.Range("K56").Value = 9147.429
$ws.Range("M56").Value = -8617.429

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H137").Value = 2385.5715
$ws.Range("I137").Value = 1799
$ws.Range("J137").Value = 2483.3333
$ws.Range("K137").Value = 5397
$ws.Range("L137").Value = 7449.999899999999
$ws.Range("N137").Value = -17649.9999
$ws.Range("M137").Value = -297

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12518.6
$ws.Range("I70").Value = 9397
$ws.Range("K70").Value = 9397
$ws.Range("M70").Value = -9127

$ws.Range("H73").Value = 12518.6
$ws.Range("I73").Value = 9397
$ws.Range("K73").Value = 9397
$ws.Range("M73").Value = -8461

$ws.Range("H80").Value = 6162.643
$ws.Range("J80").Value = 6979.625
$ws.Range("L80").Value = 6979.625
$ws.Range("N80").Value = -8975.625

$ws.Range("H83").Value = 6162.643
$ws.Range("J83").Value = 6979.625
$ws.Range("L83").Value = 34898.125
$ws.Range("N83").Value = -44882.125

$ws.Range("H102").Value = 1500.6207
$ws.Range("I102").Value = 1222.0416
$ws.Range("K102").Value = 1222.0416
$ws.Range("M102").Value = 399.9584

$ws.Range("H113").Value = 5500
$ws.Range("J113").Value = 8000
$ws.Range("L113").Value = 8000
$ws.Range("N113").Value = -12340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1862.625
$ws.Range("I68").Value = 2000.1666
$ws.Range("J68").Value = 1450
$ws.Range("K68").Value = 2000.1666
$ws.Range("L68").Value = 1450
$ws.Range("M68").Value = -1251.1666
$ws.Range("N68").Value = -2948

$ws.Range("H71").Value = 1862.625
$ws.Range("I71").Value = 2000.1666
$ws.Range("J71").Value = 1450
$ws.Range("K71").Value = 10000.833
$ws.Range("L71").Value = 7250
$ws.Range("M71").Value = -6256.833000000001
$ws.Range("N71").Value = -14738

$ws.Range("H94").Value = 100000
$ws.Range("J94").Value = 100000
$ws.Range("L94").Value = 100000
$ws.Range("N94").Value = -101352

$ws.Range("H100").Value = 1514.2
$ws.Range("I100").Value = 1514.2
$ws.Range("K100").Value = 1514.2
$ws.Range("M100").Value = -973.2

$ws.Range("H101").Value = 6262.1113
$ws.Range("J101").Value = 6262.1113
$ws.Range("L101").Value = 6262.1113
$ws.Range("N101").Value = -12752.1113

$ws.Range("H122").Value = 4203.857
$ws.Range("I122").Value = 2785.8
$ws.Range("J122").Value = 7749
$ws.Range("K122").Value = 8357.400000000001
$ws.Range("L122").Value = 23247
$ws.Range("M122").Value = -5907.400000000001
$ws.Range("N122").Value = -28147

$ws.Range("H132").Value = 4943.222
$ws.Range("I132").Value = 4415
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 13245
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -10715
$ws.Range("N132").Value = -23058.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5221.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5221.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5221.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6469.5

$ws.Range("H65").Value = 5221.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5221.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 26107.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -32347.5

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 1974.625
$ws.Range("I81").Value = 1542.1428
$ws.Range("K81").Value = 3084.2856
$ws.Range("M81").Value = -2023.2856

$ws.Range("H82").Value = 25000
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 1974.625
$ws.Range("I84").Value = 1542.1428
$ws.Range("K84").Value = 15421.428
$ws.Range("M84").Value = -10117.428

$ws.Range("H85").Value = 25000
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws.Range("H98").Value = 33999.5
$ws.Range("J98").Value = 31999.334
$ws.Range("L98").Value = 31999.334
$ws.Range("N98").Value = -37989.334

$ws.Range("H107").Value = 928.6
$ws.Range("I107").Value = 1026.4286
$ws.Range("K107").Value = 3079.2858
$ws.Range("M107").Value = -1159.2858
